$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 450.8125
$ws.Range("I2").Value = 416.64285
$ws.Range("K2").Value = 416.64285
$ws.Range("M2").Value = -303.64285
$ws.Range("H5").Value = 59.8
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H40").Value = 3216.8
$ws.Range("I40").Value = 2063
$ws.Range("J40").Value = 3711.2856
$ws.Range("K40").Value = 2063
$ws.Range("L40").Value = 3711.2856
$ws.Range("M40").Value = -1888
$ws.Range("N40").Value = -4061.2856
$ws.Range("H51").Value = 4944.0625
$ws.Range("I51").Value = 5917.4287
$ws.Range("J51").Value = 4187
$ws.Range("K51").Value = 5917.4287
$ws.Range("L51").Value = 4187
$ws.Range("M51").Value = -5433.4287
$ws.Range("N51").Value = -5155
$ws.Range("H86").Value = 25064.346
$ws.Range("I86").Value = 4545.615
$ws.Range("K86").Value = 4545.615
$ws.Range("M86").Value = -3422.615
$ws.Range("H89").Value = 25064.346
$ws.Range("I89").Value = 4545.615
$ws.Range("K89").Value = 22728.075
$ws.Range("M89").Value = -17112.075
$ws.Range("H109").Value = 205325550
$ws.Range("J109").Value = 205325550
$ws.Range("L109").Value = 205325550
$ws.Range("N109").Value = -205328324
$ws.Range("H125").Value = 6237
$ws.Range("J125").Value = 4536.75
$ws.Range("L125").Value = 40830.75
$ws.Range("N125").Value = -45750.75
$ws.Range("H127").Value = 960
$ws.Range("I127").Value = 960
$ws.Range("K127").Value = 2880
$ws.Range("M127").Value = 2080
$ws.Range("H136").Value = 69000
$ws.Range("J136").Value = 69000
$ws.Range("L136").Value = 69000
$ws.Range("N136").Value = -79200

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8073.3228
$ws.Range("I61").Value = 8219.75
$ws.Range("J61").Value = 7571.2856
$ws.Range("K61").Value = 8219.75
$ws.Range("L61").Value = 7571.2856
$ws.Range("M61").Value = -8007.75
$ws.Range("N61").Value = -7995.2856
$ws.Range("H74").Value = 4867
$ws.Range("I74").Value = 6467.1763
$ws.Range("J74").Value = 2600.0833
$ws.Range("K74").Value = 6467.1763
$ws.Range("L74").Value = 2600.0833
$ws.Range("M74").Value = -5593.1763
$ws.Range("N74").Value = -4348.0833
$ws.Range("H77").Value = 4867
$ws.Range("I77").Value = 6467.1763
$ws.Range("J77").Value = 2600.0833
$ws.Range("K77").Value = 32335.8815
$ws.Range("L77").Value = 13000.4165
$ws.Range("M77").Value = -27967.8815
$ws.Range("N77").Value = -21736.4165
$ws.Range("H132").Value = 6399.8
$ws.Range("I132").Value = 2999
$ws.Range("J132").Value = 6777.6665
$ws.Range("K132").Value = 8997
$ws.Range("L132").Value = 20332.9995
$ws.Range("M132").Value = -6467
$ws.Range("N132").Value = -25392.9995
$ws.Range("H136").Value = 8073.3228
$ws.Range("I136").Value = 8219.75
$ws.Range("J136").Value = 7571.2856
$ws.Range("K136").Value = 24659.25
$ws.Range("L136").Value = 22713.8568
$ws.Range("M136").Value = -22109.25
$ws.Range("N136").Value = -27813.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6196.75
$ws.Range("I86").Value = 6799.0625
$ws.Range("J86").Value = 3787.5
$ws.Range("K86").Value = 6799.0625
$ws.Range("L86").Value = 3787.5
$ws.Range("M86").Value = -5676.0625
$ws.Range("N86").Value = -6033.5
$ws.Range("H89").Value = 6196.75
$ws.Range("I89").Value = 6799.0625
$ws.Range("J89").Value = 3787.5
$ws.Range("K89").Value = 33995.3125
$ws.Range("L89").Value = 18937.5
$ws.Range("M89").Value = -28379.3125
$ws.Range("N89").Value = -30169.5
$ws.Range("H107").Value = 1248.2
$ws.Range("I107").Value = 1477.8572
$ws.Range("K107").Value = 1477.8572
$ws.Range("M107").Value = 442.1428000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8362.027
$ws.Range("I31").Value = 10377.368
$ws.Range("J31").Value = 6109.5884
$ws.Range("K31").Value = 10377.368
$ws.Range("L31").Value = 6109.5884
$ws.Range("M31").Value = -10082.368
$ws.Range("N31").Value = -6699.5884
$ws.Range("H34").Value = 8362.027
$ws.Range("I34").Value = 10377.368
$ws.Range("J34").Value = 6109.5884
$ws.Range("K34").Value = 10377.368
$ws.Range("L34").Value = 6109.5884
$ws.Range("M34").Value = -10175.368
$ws.Range("N34").Value = -6513.5884
$ws.Range("H62").Value = 10599.667
$ws.Range("J62").Value = 10562.5
$ws.Range("L62").Value = 10562.5
$ws.Range("N62").Value = -11810.5
$ws.Range("H65").Value = 10599.667
$ws.Range("J65").Value = 10562.5
$ws.Range("L65").Value = 52812.5
$ws.Range("N65").Value = -59052.5
$ws.Range("H74").Value = 69466.28999999999
$ws.Range("J74").Value = 69534.8
$ws.Range("L74").Value = 69534.8
$ws.Range("N74").Value = -71282.8
$ws.Range("H77").Value = 69466.28999999999
$ws.Range("J77").Value = 69534.8
$ws.Range("L77").Value = 208604.4
$ws.Range("N77").Value = -217340.4
$ws.Range("H107").Value = 6808.4443
$ws.Range("I107").Value = 8649.429
$ws.Range("K107").Value = 8649.429
$ws.Range("M107").Value = -6729.429
$ws.Range("H134").Value = 5516.5454
$ws.Range("I134").Value = 5836.12
$ws.Range("K134").Value = 17508.36
$ws.Range("M134").Value = -14973.36
$ws.Range("H137").Value = 52890
$ws.Range("J137").Value = 52890
$ws.Range("L137").Value = 52890
$ws.Range("N137").Value = -63090

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 92426056
$ws.Range("I4").Value = 88174990
$ws.Range("K4").Value = 264524970
$ws.Range("M4").Value = -264524858
$ws.Range("H11").Value = 100288.9
$ws.Range("I11").Value = 129.66667
$ws.Range("J11").Value = 143214.28
$ws.Range("K11").Value = 389.00001
$ws.Range("L11").Value = 429642.84
$ws.Range("M11").Value = -249.00001
$ws.Range("N11").Value = -429922.84
$ws.Range("H81").Value = 175096.33
$ws.Range("I81").Value = 334637.66
$ws.Range("J81").Value = 15555
$ws.Range("K81").Value = 1003912.98
$ws.Range("L81").Value = 46665
$ws.Range("M81").Value = -1002789.98
$ws.Range("N81").Value = -48911
$ws.Range("H84").Value = 175096.33
$ws.Range("I84").Value = 334637.66
$ws.Range("J84").Value = 15555
$ws.Range("K84").Value = 3011738.94
$ws.Range("L84").Value = 139995
$ws.Range("M84").Value = -3006122.94
$ws.Range("N84").Value = -151227
$ws.Range("H132").Value = 29258.955
$ws.Range("I132").Value = 825.5714
$ws.Range("K132").Value = 7430.1426
$ws.Range("M132").Value = -4900.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 28425.715
$ws.Range("J44").Value = 28425.715
$ws.Range("L44").Value = 28425.715
$ws.Range("N44").Value = -29617.715
$ws.Range("H122").Value = 12268.728
$ws.Range("I122").Value = 8688.8125
$ws.Range("K122").Value = 26066.4375
$ws.Range("M122").Value = -23616.4375
$ws.Range("H126").Value = 13837.333
$ws.Range("I126").Value = 41822.332
$ws.Range("K126").Value = 125466.996
$ws.Range("M126").Value = -122996.996
$ws.Range("H132").Value = 5313.4736
$ws.Range("I132").Value = 5313.4736
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 15940.4208
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -13410.4208
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 21373.924
$ws.Range("I7").Value = 31568.467
$ws.Range("K7").Value = 31568.467
$ws.Range("M7").Value = -31456.467
$ws.Range("H46").Value = 2193.5386
$ws.Range("J46").Value = 2751.6843
$ws.Range("L46").Value = 2751.6843
$ws.Range("N46").Value = -3127.6843
$ws.Range("H68").Value = 10727
$ws.Range("J68").Value = 10727
$ws.Range("L68").Value = 10727
$ws.Range("N68").Value = -12225
$ws.Range("H71").Value = 10727
$ws.Range("J71").Value = 10727
$ws.Range("L71").Value = 53635
$ws.Range("N71").Value = -61123
$ws.Range("H126").Value = 21373.924
$ws.Range("I126").Value = 31568.467
$ws.Range("K126").Value = 94705.401
$ws.Range("M126").Value = -92235.401
$ws.Range("H132").Value = 879616
$ws.Range("I132").Value = 2484108.8
$ws.Range("J132").Value = 4438.1816
$ws.Range("K132").Value = 7452326.399999999
$ws.Range("L132").Value = 13314.5448
$ws.Range("M132").Value = -7449796.399999999
$ws.Range("N132").Value = -18374.5448
$ws.Range("H136").Value = 4894.1177
$ws.Range("I136").Value = 3238.077
$ws.Range("K136").Value = 9714.231
$ws.Range("M136").Value = -7164.231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 15824.9
$ws.Range("J122").Value = 39094.066
$ws.Range("L122").Value = 117282.198
$ws.Range("N122").Value = -122182.198
$ws.Range("H136").Value = 364230.94
$ws.Range("I136").Value = 535979.1
$ws.Range("J136").Value = 8466.857
$ws.Range("K136").Value = 1607937.3
$ws.Range("L136").Value = 25400.571
$ws.Range("M136").Value = -1605387.3
$ws.Range("N136").Value = -30500.571

Write-Output "Applied all changes"